$d = $word.ActiveDocument

$replacements = @(
    @("2023-03-31 Friday", "2023-04-01 Saturday"),
    @("30-2=28", "89-71=18"),
    @("51+8=59", "37-32=5"),
    @("19+39=58", "95-12=83"),
    @("21-18=3", "1+83=84"),
    @("68+13=81", "91-72=19"),
    @("19+35=54", "63-32=31"),
    @("39+19=58", "38+27=65"),
    @("67-38=29", "25-19=6"),
    @("99-19=80", "47+51=98"),
    @("83-20=63", "88+8=96"),
    @("60+7=67", "31+25=56"),
    @("89-79=10", "23+43=66"),
    @("62+34=96", "0+39=39"),
    @("18+51=69", "67+9=76"),
    @("23+35=58", "68-0=68"),
    @("53-10=43", "1+50=51"),
    @("54-40=14", "59-41=18"),
    @("87+11=98", "76-47=29"),
    @("26+59=85", "43-22=21"),
    @("72-25=47", "83-10=73"),
    @("82-29=53", "86-74=12"),
    @("10+52=62", "84-67=17"),
    @("47+0=47", "90-41=49"),
    @("20+48=68", "36+44=80"),
    @("33-14=19", "23+62=85"),
    @("23-22=1", "46-14=32"),
    @("31-21=10", "78-18=60"),
    @("37-36=1", "23+74=97"),
    @("10+89=99", "81-39=42"),
    @("69+26=95", "25+74=99"),
    @("78-52=26", "44-6=38"),
    @("51+11=62", "98-68=30"),
    @("98-59=39", "70+15=85"),
    @("80-31=49", "56-20=36"),
    @("98+1=99", "61-28=33"),
    @("62-27=35", "53+37=90"),
    @("67-26=41", "37-12=25"),
    @("37-22=15", "7+45=52"),
    @("68+23=91", "34+23=57"),
    @("9+44=53", "25+47=72"),
    @("68+24=92", "25+72=97"),
    @("26+43=69", "22+60=82"),
    @("0+88=88", "47-1=46"),
    @("12+58=70", "80-44=36"),
    @("75+24=99", "31+12=43"),
    @("61-40=21", "43+55=98"),
    @("44+3=47", "66-25=41"),
    @("72-6=66", "41-15=26"),
    @("25+0=25", "21-12=9"),
    @("36-14=22", "27-21=6"),
    @("53-24=29", "41+5=46"),
    @("67-41=26", "37-12=25"),
    @("75-75=0", "85-83=2"),
    @("22+14=36", "3+10=13"),
    @("19+49=68", "81-19=62"),
    @("75-11=64", "59+10=69"),
    @("16+30=46", "95-19=76"),
    @("1+32=33", "79+18=97"),
    @("30+34=64", "11+5=16"),
    @("78-16=62", "6+46=52"),
    @("14-4=10", "21+13=34"),
    @("33-5=28", "67-56=11"),
    @("66-6=60", "0+65=65"),
    @("3+68=71", "1+90=91"),
    @("66+0=66", "77-15=62"),
    @("25+55=80", "67-18=49"),
    @("49+27=76", "77-14=63"),
    @("30+11=41", "95-6=89"),
    @("8+91=99", "79-50=29"),
    @("65-5=60", "30+48=78"),
    @("58-1=57", "60+39=99"),
    @("72-42=30", "62+32=94"),
    @("43-30=13", "97-54=43"),
    @("49-4=45", "56-1=55"),
    @("99-97=2", "60+6=66"),
    @("33+63=96", "9+40=49"),
    @("74+23=97", "96-90=6"),
    @("18+53=71", "76+3=79"),
    @("57+21=78", "87-33=54"),
    @("84-65=19", "90-73=17"),
    @("44-43=1", "36-22=14"),
    @("29-14=15", "27-5=22"),
    @("3+23=26", "1+58=59"),
    @("4+32=36", "29+24=53"),
    @("23+27=50", "54-27=27"),
    @("59-9=50", "43+45=88"),
    @("83-42=41", "91-63=28"),
    @("73-66=7", "76-62=14"),
    @("0+64=64", "62-0=62"),
    @("72-3=69", "96-70=26"),
    @("50+39=89", "30+30=60"),
    @("46+41=87", "10+62=72"),
    @("64-61=3", "70-39=31"),
    @("38-3=35", "9+10=19"),
    @("68-32=36", "91-58=33"),
    @("54-35=19", "11+12=23"),
    @("19-11=8", "27-10=17"),
    @("14+58=72", "27-10=17"),
    @("52+32=84", "69-61=8"),
    @("97-83=14", "74-71=3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($replacements.Count) text runs."
